$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values look numeric,
# so Excel keeps them as text (matching original inlineStr text cells)
# instead of converting to a number (which would drop formatting like
# trailing zeros, e.g. "1.00" -> 1).
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated values
$ws.Range('D2').Value = '51.089.88'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.958.03'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '380.48'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').Value = '102.07'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  +1.70%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.590'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').Value = '36.41'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.423.22'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '7.78'
$ws.Range('E14').Value = '  +5.32%  '
$ws.Range('D15').Value = '18.31'
$ws.Range('E15').Value = '  +2.38%  '
$ws.Range('D16').Value = '11.26'
$ws.Range('E16').Value = '  +7.06%  '
$ws.Range('D17').Value = '2.951.96'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '0.996'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '51.180.29'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').Value = '12.36'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').Value = '3.32'
$ws.Range('E23').Value = '  +11.01%  '
$ws.Range('D24').Value = '70.25'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('D25').Value = '266.90'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').Value = '7.87'
$ws.Range('E26').Value = '  -5.25%  '
$ws.Range('D27').Value = '7.19'
$ws.Range('E27').Value = '  -9.36%  '
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '0.165'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D30').Value = '25.83'
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = '10.28'
$ws.Range('E32').Value = '  +4.32%  '
$ws.Range('D33').Value = '34.32'
$ws.Range('E33').Value = '  +2.27%  '
$ws.Range('D34').Value = '51.06'
$ws.Range('E34').Value = '  +0.76%  '
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').Value = '0.0435'
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = '3.22'
$ws.Range('E38').Value = '  +7.54%  '
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +0.82%  '
$ws.Range('E40').Value = '  +2.26%  '
$ws.Range('D41').Value = '16.49'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '2.50'
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('D43').Value = '124.72'
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('D44').Value = '3.52'
$ws.Range('E44').Value = '  +7.05%  '
$ws.Range('D45').Value = '21.34'
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('D46').Value = '0.273'
$ws.Range('E46').Value = '  -6.85%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  +2.84%  '
$ws.Range('D49').Value = '2.045.24'
$ws.Range('E49').Value = '  +2.96%  '
$ws.Range('D50').Value = '0.0320'
$ws.Range('E50').Value = '  -6.30%  '
$ws.Range('D51').Value = '5.40'
$ws.Range('E51').Value = '  +6.94%  '
